$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add headers E1/F1 on sheet1, copying header style from D1 ---
$ws1.Cells.Item(1,5).Value = "File size"
$ws1.Cells.Item(1,6).Value = "More info"
$ws1.Range("D1").Copy()
$ws1.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update D column (URL) with inserted-space variant + add E (File size) and F (More info) ---
$ws1.Cells.Item(2,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/5935EPStudies/eps dtl32139.pdf '
$ws1.Cells.Item(2,5).Value = '388 kb '
$ws1.Cells.Item(2,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= 5935EPStudies '
$ws1.Cells.Item(3,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-DTL-17/dtl17s s224.pdf '
$ws1.Cells.Item(3,5).Value = '117 kb '
$ws1.Cells.Item(3,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-DTL-17 '
$ws1.Cells.Item(4,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-DTL-39030/dtl 39030ss20.pdf '
$ws1.Cells.Item(4,5).Value = '130 kb '
$ws1.Cells.Item(4,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-DTL-39030 '
$ws1.Cells.Item(5,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-DTL-39030/dtl 39030ss21.pdf '
$ws1.Cells.Item(5,5).Value = '108 kb '
$ws1.Cells.Item(5,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-DTL-39030 '
$ws1.Cells.Item(6,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-DTL-83503/dtl 83503.pdf '
$ws1.Cells.Item(6,5).Value = '288 kb '
$ws1.Cells.Item(6,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-DTL-83503 '
$ws1.Cells.Item(7,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-1/prf1ss3 0not1.pdf '
$ws1.Cells.Item(7,5).Value = '61 kb '
$ws1.Cells.Item(7,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-1 '
$ws1.Cells.Item(8,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-1/prf1ss1 48not1.pdf '
$ws1.Cells.Item(8,5).Value = '60 kb '
$ws1.Cells.Item(8,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-1 '
$ws1.Cells.Item(9,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-1/prf1ss7 81not1.pdf '
$ws1.Cells.Item(9,5).Value = '60 kb '
$ws1.Cells.Item(9,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-1 '
$ws1.Cells.Item(10,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-1/prf1ss1 047not1.pdf '
$ws1.Cells.Item(10,5).Value = '61 kb '
$ws1.Cells.Item(10,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-1 '
$ws1.Cells.Item(11,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-1/prf1ss1 636.pdf '
$ws1.Cells.Item(11,5).Value = '278 kb '
$ws1.Cells.Item(11,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-1 '
$ws1.Cells.Item(12,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-19500/idp rf19500ss426.pdf '
$ws1.Cells.Item(12,5).Value = '549 kb '
$ws1.Cells.Item(12,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-19500 '
$ws1.Cells.Item(13,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-19500/idp rf19500ss439.pdf '
$ws1.Cells.Item(13,5).Value = '637 kb '
$ws1.Cells.Item(13,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-19500 '
$ws1.Cells.Item(14,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-19500/idp rf19500ss782.pdf '
$ws1.Cells.Item(14,5).Value = '1007 kb '
$ws1.Cells.Item(14,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-19500 '
$ws1.Cells.Item(15,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-39016/prf 39016ss48.pdf '
$ws1.Cells.Item(15,5).Value = '160 kb '
$ws1.Cells.Item(15,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-39016 '
$ws1.Cells.Item(16,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-55339/prf 55339ss39.pdf '
$ws1.Cells.Item(16,5).Value = '244 kb '
$ws1.Cells.Item(16,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-55339 '
$ws1.Cells.Item(17,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-55339/prf 55339ss48.pdf '
$ws1.Cells.Item(17,5).Value = '116 kb '
$ws1.Cells.Item(17,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-55339 '
$ws1.Cells.Item(18,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-55339/prf 55339ss51.pdf '
$ws1.Cells.Item(18,5).Value = '196 kb '
$ws1.Cells.Item(18,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-55339 '
$ws1.Cells.Item(19,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-55339/prf 55339ss54.pdf '
$ws1.Cells.Item(19,5).Value = '121 kb '
$ws1.Cells.Item(19,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-55339 '
$ws1.Cells.Item(20,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-55339/prf 55339ss55.pdf '
$ws1.Cells.Item(20,5).Value = '137 kb '
$ws1.Cells.Item(20,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-55339 '
$ws1.Cells.Item(21,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-6106/idpr f6106.pdf '
$ws1.Cells.Item(21,5).Value = '862 kb '
$ws1.Cells.Item(21,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-6106 '
$ws1.Cells.Item(22,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-6106/idpr f6106sup1.pdf '
$ws1.Cells.Item(22,5).Value = '192 kb '
$ws1.Cells.Item(22,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-6106 '
$ws1.Cells.Item(23,4).Value = 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-6106/ms27 706.pdf '
$ws1.Cells.Item(23,5).Value = '381 kb '
$ws1.Cells.Item(23,6).Value = 'https://landandmaritimeapps.dla.mil/Programs/MilSpec/ListDocs.aspx?BasicDoc= MIL-PRF-6106 '
# --- Create sheet2 ("pagina1") as a copy of sheet1, then trim to 6 rows ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "pagina1"
$ws2.Range("A7:F23").Clear()
$ws2.Cells.Item(6,6).Value = ""
$ws2.Cells.Item(6,6).ClearContents()

# --- Create sheet3 ("pagina2") as a copy of sheet2 (pagina1), then tweak F6 to an empty (but present) cell ---
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "pagina2"
$ws3.Cells.Item(6,6).Font.Bold = $false

# Restore original active sheet/tab selection
$ws1.Activate()
